$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1039788.7
$ws.Cells.Item(2, 9).Value = 1454903.8
$ws.Cells.Item(2, 10).Value = 2001
$ws.Cells.Item(2, 11).Value = 1454903.8
$ws.Cells.Item(2, 12).Value = 2001
$ws.Cells.Item(2, 13).Value = -1454790.8
$ws.Cells.Item(2, 14).Value = -2227

$ws.Cells.Item(6, 8).Value = 1363.8334
$ws.Cells.Item(6, 9).Value = 545.75
$ws.Cells.Item(6, 10).Value = 3000
$ws.Cells.Item(6, 11).Value = 1637.25
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = -1525.25
$ws.Cells.Item(6, 14).Value = -9224

$ws.Cells.Item(9, 8).Value = 1116.6666
$ws.Cells.Item(9, 9).Value = 424
$ws.Cells.Item(9, 10).Value = 2502
$ws.Cells.Item(9, 11).Value = 424
$ws.Cells.Item(9, 12).Value = 2502
$ws.Cells.Item(9, 13).Value = -255
$ws.Cells.Item(9, 14).Value = -2840

$ws.Cells.Item(12, 8).Value = 525
$ws.Cells.Item(12, 9).Value = 366.66666
$ws.Cells.Item(12, 10).Value = 1000
$ws.Cells.Item(12, 11).Value = 366.66666
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 13).Value = -196.66666
$ws.Cells.Item(12, 14).Value = -1340

$ws.Cells.Item(21, 8).Value = 22899.6
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 22899.6
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 22899.6
$ws.Cells.Item(21, 14).Value = -23835.6

$ws.Cells.Item(23, 8).Value = 22899.6
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 22899.6
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 22899.6
$ws.Cells.Item(23, 14).Value = -23367.6

$ws.Cells.Item(29, 8).Value = 4933.3335
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 4933.3335
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 14800.0005
$ws.Cells.Item(29, 14).Value = -15362.0005

$ws.Cells.Item(58, 8).Value = 4324.385
$ws.Cells.Item(58, 9).Value = 316.85715
$ws.Cells.Item(58, 10).Value = 8999.833000000001
$ws.Cells.Item(58, 11).Value = 950.5714499999999
$ws.Cells.Item(58, 12).Value = 26999.499
$ws.Cells.Item(58, 13).Value = -800.5714499999999
$ws.Cells.Item(58, 14).Value = -27299.499

$ws.Cells.Item(96, 8).Value = 1742.7778
$ws.Cells.Item(96, 9).Value = 1197.8334
$ws.Cells.Item(96, 10).Value = 2832.6667
$ws.Cells.Item(96, 11).Value = 3593.5002
$ws.Cells.Item(96, 12).Value = 8498.000100000001
$ws.Cells.Item(96, 13).Value = -2220.5002
$ws.Cells.Item(96, 14).Value = -11244.0001

$ws.Cells.Item(106, 8).Value = 6956.8
$ws.Cells.Item(106, 9).Value = 4445.25
$ws.Cells.Item(106, 10).Value = 17003
$ws.Cells.Item(106, 11).Value = 4445.25
$ws.Cells.Item(106, 12).Value = 17003
$ws.Cells.Item(106, 13).Value = -3814.25
$ws.Cells.Item(106, 14).Value = -18265

$ws.Cells.Item(116, 8).Value = 3355.7144
$ws.Cells.Item(116, 9).Value = 3332.3333
$ws.Cells.Item(116, 10).Value = 3496
$ws.Cells.Item(116, 11).Value = 3332.3333
$ws.Cells.Item(116, 12).Value = 3496
$ws.Cells.Item(116, 13).Value = 109.6667000000002
$ws.Cells.Item(116, 14).Value = -10380

$ws.Cells.Item(132, 8).Value = 4231.846
$ws.Cells.Item(132, 9).Value = 3554
$ws.Cells.Item(132, 10).Value = 29990
$ws.Cells.Item(132, 11).Value = 10662
$ws.Cells.Item(132, 12).Value = 89970
$ws.Cells.Item(132, 13).Value = -8132
$ws.Cells.Item(132, 14).Value = -95030

$ws.Cells.Item(137, 8).Value = 927603.7
$ws.Cells.Item(137, 9).Value = 1315
$ws.Cells.Item(137, 10).Value = 3706469.8
$ws.Cells.Item(137, 11).Value = 3945
$ws.Cells.Item(137, 12).Value = 11119409.4
$ws.Cells.Item(137, 13).Value = -1395
$ws.Cells.Item(137, 14).Value = -11124509.4

$ws.Cells.Item(141, 9).Value = 2246.2666
$ws.Cells.Item(141, 10).Value = 1445
$ws.Cells.Item(141, 11).Value = 6738.7998
$ws.Cells.Item(141, 12).Value = 4335
$ws.Cells.Item(141, 13).Value = -1558.7998
$ws.Cells.Item(141, 14).Value = -14695

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17017.61
$ws.Cells.Item(32, 9).Value = 19849.754
$ws.Cells.Item(32, 10).Value = 3371.818
$ws.Cells.Item(32, 11).Value = 19849.754
$ws.Cells.Item(32, 12).Value = 3371.818
$ws.Cells.Item(32, 13).Value = -19562.754
$ws.Cells.Item(32, 14).Value = -3945.818

$ws.Cells.Item(43, 8).Value = 20000
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 20000
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 20000
$ws.Cells.Item(43, 14).Value = -20626

$ws.Cells.Item(45, 8).Value = 2564.3157
$ws.Cells.Item(45, 9).Value = 1763.9
$ws.Cells.Item(45, 10).Value = 3453.6667
$ws.Cells.Item(45, 11).Value = 1763.9
$ws.Cells.Item(45, 12).Value = 3453.6667
$ws.Cells.Item(45, 13).Value = -1386.9
$ws.Cells.Item(45, 14).Value = -4207.6667

$ws.Cells.Item(74, 8).Value = 1549.3658
$ws.Cells.Item(74, 9).Value = 1132.1143
$ws.Cells.Item(74, 10).Value = 3983.3333
$ws.Cells.Item(74, 11).Value = 1132.1143
$ws.Cells.Item(74, 12).Value = 3983.3333
$ws.Cells.Item(74, 13).Value = -258.1143
$ws.Cells.Item(74, 14).Value = -5731.3333

$ws.Cells.Item(77, 8).Value = 1549.3658
$ws.Cells.Item(77, 9).Value = 1132.1143
$ws.Cells.Item(77, 10).Value = 3983.3333
$ws.Cells.Item(77, 11).Value = 5660.5715
$ws.Cells.Item(77, 12).Value = 19916.6665
$ws.Cells.Item(77, 13).Value = -1292.5715
$ws.Cells.Item(77, 14).Value = -28652.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2499.5
$ws.Cells.Item(20, 9).Value = 1070.8572
$ws.Cells.Item(20, 10).Value = 4499.6
$ws.Cells.Item(20, 11).Value = 1070.8572
$ws.Cells.Item(20, 12).Value = 4499.6
$ws.Cells.Item(20, 13).Value = -823.8571999999999
$ws.Cells.Item(20, 14).Value = -4993.6

$ws.Cells.Item(102, 8).Value = 6215.25
$ws.Cells.Item(102, 9).Value = 6215.25
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 6215.25
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -2970.25

$ws.Cells.Item(134, 8).Value = 1811.8148
$ws.Cells.Item(134, 9).Value = 1814.56
$ws.Cells.Item(134, 10).Value = 1777.5
$ws.Cells.Item(134, 11).Value = 5443.68
$ws.Cells.Item(134, 12).Value = 5332.5
$ws.Cells.Item(134, 13).Value = -2908.68
$ws.Cells.Item(134, 14).Value = -10402.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1631.6316
$ws.Cells.Item(31, 9).Value = 1611.1666
$ws.Cells.Item(31, 10).Value = 2000
$ws.Cells.Item(31, 11).Value = 1611.1666
$ws.Cells.Item(31, 12).Value = 2000
$ws.Cells.Item(31, 13).Value = -1316.1666
$ws.Cells.Item(31, 14).Value = -2590

$ws.Cells.Item(34, 8).Value = 1631.6316
$ws.Cells.Item(34, 9).Value = 1611.1666
$ws.Cells.Item(34, 10).Value = 2000
$ws.Cells.Item(34, 11).Value = 1611.1666
$ws.Cells.Item(34, 12).Value = 2000
$ws.Cells.Item(34, 13).Value = -1409.1666
$ws.Cells.Item(34, 14).Value = -2404

$ws.Cells.Item(107, 8).Value = 2114.258
$ws.Cells.Item(107, 9).Value = 208.35294
$ws.Cells.Item(107, 10).Value = 4428.5713
$ws.Cells.Item(107, 11).Value = 208.35294
$ws.Cells.Item(107, 12).Value = 4428.5713
$ws.Cells.Item(107, 13).Value = 1711.64706
$ws.Cells.Item(107, 14).Value = -8268.5713

$ws.Cells.Item(132, 8).Value = 2010.95
$ws.Cells.Item(132, 9).Value = 1870.5
$ws.Cells.Item(132, 10).Value = 3275
$ws.Cells.Item(132, 11).Value = 5611.5
$ws.Cells.Item(132, 12).Value = 9825
$ws.Cells.Item(132, 13).Value = -3081.5
$ws.Cells.Item(132, 14).Value = -14885

$ws.Cells.Item(134, 8).Value = 68696.60000000001
$ws.Cells.Item(134, 9).Value = 73503.42999999999
$ws.Cells.Item(134, 10).Value = 1401
$ws.Cells.Item(134, 11).Value = 220510.29
$ws.Cells.Item(134, 12).Value = 4203
$ws.Cells.Item(134, 13).Value = -217975.29
$ws.Cells.Item(134, 14).Value = -9273

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 668.7273
$ws.Cells.Item(5, 9).Value = 375.1111
$ws.Cells.Item(5, 10).Value = 1990
$ws.Cells.Item(5, 11).Value = 1125.3333
$ws.Cells.Item(5, 12).Value = 5970
$ws.Cells.Item(5, 13).Value = -1013.3333
$ws.Cells.Item(5, 14).Value = -6194

$ws.Cells.Item(68, 8).Value = 39873.23
$ws.Cells.Item(68, 9).Value = 1474
$ws.Cells.Item(68, 10).Value = 46854.91
$ws.Cells.Item(68, 11).Value = 4422
$ws.Cells.Item(68, 12).Value = 140564.73
$ws.Cells.Item(68, 13).Value = -3611
$ws.Cells.Item(68, 14).Value = -142186.73

$ws.Cells.Item(71, 8).Value = 39873.23
$ws.Cells.Item(71, 9).Value = 1474
$ws.Cells.Item(71, 10).Value = 46854.91
$ws.Cells.Item(71, 11).Value = 13266
$ws.Cells.Item(71, 12).Value = 421694.1900000001
$ws.Cells.Item(71, 13).Value = -9210
$ws.Cells.Item(71, 14).Value = -429806.1900000001

$ws.Cells.Item(127, 8).Value = 7781.6665
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 7781.6665
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 23344.9995
$ws.Cells.Item(127, 14).Value = -33264.99950000001

$ws.Cells.Item(135, 8).Value = 668.7273
$ws.Cells.Item(135, 9).Value = 375.1111
$ws.Cells.Item(135, 10).Value = 1990
$ws.Cells.Item(135, 11).Value = 3375.9999
$ws.Cells.Item(135, 12).Value = 17910
$ws.Cells.Item(135, 13).Value = -840.9999000000003
$ws.Cells.Item(135, 14).Value = -22980

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6504
$ws.Cells.Item(70, 9).Value = 6504
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 6504
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = -6234

$ws.Cells.Item(73, 8).Value = 6504
$ws.Cells.Item(73, 9).Value = 6504
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 6504
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = -5568

$ws.Cells.Item(102, 8).Value = 3954
$ws.Cells.Item(102, 9).Value = 5506
$ws.Cells.Item(102, 10).Value = 850
$ws.Cells.Item(102, 11).Value = 5506
$ws.Cells.Item(102, 12).Value = 850
$ws.Cells.Item(102, 13).Value = -3884
$ws.Cells.Item(102, 14).Value = -4094

$ws.Cells.Item(132, 8).Value = 37802.17
$ws.Cells.Item(132, 9).Value = 43444.96
$ws.Cells.Item(132, 10).Value = 2534.75
$ws.Cells.Item(132, 11).Value = 130334.88
$ws.Cells.Item(132, 12).Value = 7604.25
$ws.Cells.Item(132, 13).Value = -127804.88
$ws.Cells.Item(132, 14).Value = -12664.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5171.615
$ws.Cells.Item(7, 9).Value = 4466
$ws.Cells.Item(7, 10).Value = 5383.3
$ws.Cells.Item(7, 11).Value = 4466
$ws.Cells.Item(7, 12).Value = 5383.3
$ws.Cells.Item(7, 13).Value = -4354
$ws.Cells.Item(7, 14).Value = -5607.3

$ws.Cells.Item(61, 8).Value = 6358.7383
$ws.Cells.Item(61, 9).Value = 5696.2188
$ws.Cells.Item(61, 10).Value = 8478.799999999999
$ws.Cells.Item(61, 11).Value = 5696.2188
$ws.Cells.Item(61, 12).Value = 8478.799999999999
$ws.Cells.Item(61, 13).Value = -5494.2188
$ws.Cells.Item(61, 14).Value = -8882.799999999999

$ws.Cells.Item(82, 8).Value = 2980.3333
$ws.Cells.Item(82, 9).Value = 882
$ws.Cells.Item(82, 10).Value = 3400
$ws.Cells.Item(82, 11).Value = 882
$ws.Cells.Item(82, 12).Value = 3400
$ws.Cells.Item(82, 13).Value = -521
$ws.Cells.Item(82, 14).Value = -4122

$ws.Cells.Item(85, 8).Value = 2980.3333
$ws.Cells.Item(85, 9).Value = 882
$ws.Cells.Item(85, 10).Value = 3400
$ws.Cells.Item(85, 11).Value = 882
$ws.Cells.Item(85, 12).Value = 3400
$ws.Cells.Item(85, 13).Value = 366
$ws.Cells.Item(85, 14).Value = -5896

$ws.Cells.Item(113, 8).Value = 6358.7383
$ws.Cells.Item(113, 9).Value = 5696.2188
$ws.Cells.Item(113, 10).Value = 8478.799999999999
$ws.Cells.Item(113, 11).Value = 5696.2188
$ws.Cells.Item(113, 12).Value = 8478.799999999999
$ws.Cells.Item(113, 13).Value = -3526.2188
$ws.Cells.Item(113, 14).Value = -12818.8

$ws.Cells.Item(126, 8).Value = 5171.615
$ws.Cells.Item(126, 9).Value = 4466
$ws.Cells.Item(126, 10).Value = 5383.3
$ws.Cells.Item(126, 11).Value = 13398
$ws.Cells.Item(126, 12).Value = 16149.9
$ws.Cells.Item(126, 13).Value = -10928
$ws.Cells.Item(126, 14).Value = -21089.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2277.625
$ws.Cells.Item(81, 9).Value = 2360.1428
$ws.Cells.Item(81, 10).Value = 1700
$ws.Cells.Item(81, 11).Value = 4720.2856
$ws.Cells.Item(81, 12).Value = 3400
$ws.Cells.Item(81, 13).Value = -3659.2856
$ws.Cells.Item(81, 14).Value = -5522

$ws.Cells.Item(84, 8).Value = 2277.625
$ws.Cells.Item(84, 9).Value = 2360.1428
$ws.Cells.Item(84, 10).Value = 1700
$ws.Cells.Item(84, 11).Value = 23601.428
$ws.Cells.Item(84, 12).Value = 17000
$ws.Cells.Item(84, 13).Value = -18297.428
$ws.Cells.Item(84, 14).Value = -27608

$ws.Cells.Item(126, 8).Value = 4854.1665
$ws.Cells.Item(126, 9).Value = 4675.2144
$ws.Cells.Item(126, 10).Value = 5480.5
$ws.Cells.Item(126, 11).Value = 14025.6432
$ws.Cells.Item(126, 12).Value = 16441.5
$ws.Cells.Item(126, 13).Value = -11555.6432
$ws.Cells.Item(126, 14).Value = -21381.5

$ws.Cells.Item(135, 8).Value = 90979.57000000001
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 90979.57000000001
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 90979.57000000001
$ws.Cells.Item(135, 14).Value = -101119.57

$ws.Cells.Item(136, 8).Value = 3185.1667
$ws.Cells.Item(136, 9).Value = 2323.5715
$ws.Cells.Item(136, 10).Value = 6200.75
$ws.Cells.Item(136, 11).Value = 6970.7145
$ws.Cells.Item(136, 12).Value = 18602.25
$ws.Cells.Item(136, 13).Value = -4420.7145
$ws.Cells.Item(136, 14).Value = -23702.25
